# Fix typos in test report (test-report.xlsx)
#
# This script reproduces the textual corrections described in the commit
# "Fix typos in test report":
#   - "Passedword"                         -> "Password"
#   - "Succesfull"                         -> "Successful"
#   - "Rate dish ..."                      -> "Rate dishes ..."
#   - rich-text description in F60 reworded + split into two runs

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Keep gridlines visible (matches the original authoring state).
$excel.ActiveWindow.DisplayGridlines = $true

# --- Whole-sheet text fixes -------------------------------------------------
# Use Range.Replace so every cell sharing the mistyped string gets fixed at
# once (several rows reuse the exact same shared string).
$all = $ws.Range("A1:H89")

$all.Replace("Passedword", "Password")
$all.Replace("Succesfull", "Successful")
$all.Replace("Rate dish ", "Rate dishes ")

# --- F60: reworded test case description, split across two runs -----------
$f60 = $ws.Range("F60")
$f60.Value = "Show list of dishes which has a keyword. If there is no such a dish, the return blank result"

# Give the trailing word "result" its own run (mirrors the source diff,
# which keeps the sentence + the final word "result" as distinct runs).
$tail = $f60.Characters(87, 6)
$tail.Font.Name = "Arial"
$tail.Font.Size = 11
$tail.Font.Color = 0

# --- Restore selection / scroll position -----------------------------------
$ws.Range("J17").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "Typo fixes applied"
